$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 23558.387
$ws.Range("J17").Value = 23558.387
$ws.Range("L17").Value = 70675.16099999999
$ws.Range("N17").Value = -71011.16099999999

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 108.70588
$ws.Range("I33").Value = 96.96774000000001
$ws.Range("J33").Value = 230
$ws.Range("K33").Value = 96.96774000000001
$ws.Range("L33").Value = 230
$ws.Range("M33").Value = 132.03226
$ws.Range("N33").Value = -688

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3268.8223
$ws.Range("I64").Value = 3069
$ws.Range("J64").Value = 3350
$ws.Range("K64").Value = 3069
$ws.Range("L64").Value = 3350
$ws.Range("M64").Value = -2821
$ws.Range("N64").Value = -3846

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3268.8223
$ws.Range("I67").Value = 3069
$ws.Range("J67").Value = 3350
$ws.Range("K67").Value = 3069
$ws.Range("L67").Value = 3350
$ws.Range("M67").Value = -2211
$ws.Range("N67").Value = -5066

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 13335680
$ws.Range("I116").Value = 33335134
$ws.Range("J116").Value = 2711.6667
$ws.Range("K116").Value = 33335134
$ws.Range("L116").Value = 2711.6667
$ws.Range("M116").Value = -33331692
$ws.Range("N116").Value = -9595.6667

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1076.8718
$ws.Range("I129").Value = 733.3333
$ws.Range("J129").Value = 1105.5
$ws.Range("K129").Value = 2199.9999
$ws.Range("L129").Value = 3316.5
$ws.Range("M129").Value = 2800.0001
$ws.Range("N129").Value = -13316.5

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 678.0769
$ws.Range("I135").Value = 571.4894
$ws.Range("J135").Value = 1680
$ws.Range("K135").Value = 5143.404600000001
$ws.Range("L135").Value = 15120
$ws.Range("M135").Value = -2608.404600000001
$ws.Range("N135").Value = -20190

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1269
$ws.Range("I137").Value = 1134.7715
$ws.Range("J137").Value = 1545.3529
$ws.Range("K137").Value = 3404.3145
$ws.Range("L137").Value = 4636.0587
$ws.Range("M137").Value = -854.3145000000004
$ws.Range("N137").Value = -9736.058700000001

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3613.804
$ws.Range("I138").Value = 2817.037
$ws.Range("J138").Value = 4510.1665
$ws.Range("K138").Value = 8451.110999999999
$ws.Range("L138").Value = 13530.4995
$ws.Range("M138").Value = -3311.110999999999
$ws.Range("N138").Value = -23810.4995

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5178.543
$ws.Range("I141").Value = 2313.875
$ws.Range("K141").Value = 6941.625
$ws.Range("M141").Value = -1761.625

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9535.187
$ws.Range("I32").Value = 10755.362
$ws.Range("J32").Value = 5372.2354
$ws.Range("K32").Value = 10755.362
$ws.Range("L32").Value = 5372.2354
$ws.Range("M32").Value = -10468.362
$ws.Range("N32").Value = -5946.2354

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1165.6154
$ws.Range("I61").Value = 1109.4762
$ws.Range("J61").Value = 1401.4
$ws.Range("K61").Value = 1109.4762
$ws.Range("L61").Value = 1401.4
$ws.Range("M61").Value = -897.4762000000001
$ws.Range("N61").Value = -1825.4

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 737.5349
$ws.Range("I74").Value = 671.79486
$ws.Range("J74").Value = 1378.5
$ws.Range("K74").Value = 671.79486
$ws.Range("L74").Value = 1378.5
$ws.Range("M74").Value = 202.20514
$ws.Range("N74").Value = -3126.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 737.5349
$ws.Range("I77").Value = 671.79486
$ws.Range("J77").Value = 1378.5
$ws.Range("K77").Value = 3358.9743
$ws.Range("L77").Value = 6892.5
$ws.Range("M77").Value = 1009.0257
$ws.Range("N77").Value = -15628.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1039.619
$ws.Range("I122").Value = 981
$ws.Range("J122").Value = 1134.875
$ws.Range("K122").Value = 2943
$ws.Range("L122").Value = 3404.625
$ws.Range("M122").Value = -493
$ws.Range("N122").Value = -8304.625

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2210.543
$ws.Range("I132").Value = 1545.5294
$ws.Range("J132").Value = 2838.611
$ws.Range("K132").Value = 4636.5882
$ws.Range("L132").Value = 8515.832999999999
$ws.Range("M132").Value = -2106.5882
$ws.Range("N132").Value = -13575.833

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1165.6154
$ws.Range("I136").Value = 1109.4762
$ws.Range("J136").Value = 1401.4
$ws.Range("K136").Value = 3328.4286
$ws.Range("L136").Value = 4204.200000000001
$ws.Range("M136").Value = -778.4286000000002
$ws.Range("N136").Value = -9304.200000000001

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 127462.5
$ws.Range("I94").Value = 1925
$ws.Range("J94").Value = 253000
$ws.Range("K94").Value = 1925
$ws.Range("L94").Value = 253000
$ws.Range("M94").Value = -1474
$ws.Range("N94").Value = -253902

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2623.3872
$ws.Range("I134").Value = 2257
$ws.Range("J134").Value = 4150
$ws.Range("K134").Value = 6771
$ws.Range("L134").Value = 12450
$ws.Range("M134").Value = -4236
$ws.Range("N134").Value = -17520

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1950.7291
$ws.Range("I31").Value = 1377.258
$ws.Range("J31").Value = 2996.4707
$ws.Range("K31").Value = 1377.258
$ws.Range("L31").Value = 2996.4707
$ws.Range("M31").Value = -1082.258
$ws.Range("N31").Value = -3586.4707

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1950.7291
$ws.Range("I34").Value = 1377.258
$ws.Range("J34").Value = 2996.4707
$ws.Range("K34").Value = 1377.258
$ws.Range("L34").Value = 2996.4707
$ws.Range("M34").Value = -1175.258
$ws.Range("N34").Value = -3400.4707

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1685641.6
$ws.Range("I58").Value = 4117239
$ws.Range("J58").Value = 2227.923
$ws.Range("K58").Value = 4117239
$ws.Range("L58").Value = 2227.923
$ws.Range("M58").Value = -4117036
$ws.Range("N58").Value = -2633.923

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 565226.3
$ws.Range("I132").Value = 712626.5600000001
$ws.Range("J132").Value = 5105.2
$ws.Range("K132").Value = 2137879.68
$ws.Range("L132").Value = 15315.6
$ws.Range("M132").Value = -2135349.68
$ws.Range("N132").Value = -20375.6

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1685641.6
$ws.Range("I136").Value = 4117239
$ws.Range("J136").Value = 2227.923
$ws.Range("K136").Value = 12351717
$ws.Range("L136").Value = 6683.768999999999
$ws.Range("M136").Value = -12349167
$ws.Range("N136").Value = -11783.769

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2225
$ws.Range("I138").Value = 993.8
$ws.Range("J138").Value = 3251
$ws.Range("K138").Value = 2981.4
$ws.Range("L138").Value = 9753
$ws.Range("M138").Value = 2158.6
$ws.Range("N138").Value = -20033

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3597.3572
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3279.3572
$ws.Range("I80").Value = 3511.6667
$ws.Range("J80").Value = 2861.2
$ws.Range("K80").Value = 3511.6667
$ws.Range("L80").Value = 2861.2
$ws.Range("M80").Value = -2513.6667
$ws.Range("N80").Value = -4857.2

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3279.3572
$ws.Range("I83").Value = 3511.6667
$ws.Range("J83").Value = 2861.2
$ws.Range("K83").Value = 17558.3335
$ws.Range("L83").Value = 14306
$ws.Range("M83").Value = -12566.3335
$ws.Range("N83").Value = -24290

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2217.5
$ws.Range("I126").Value = 1545.1428
$ws.Range("J126").Value = 3158.8
$ws.Range("K126").Value = 4635.428400000001
$ws.Range("L126").Value = 9476.400000000001
$ws.Range("M126").Value = -2165.428400000001
$ws.Range("N126").Value = -14416.4

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2486.2942
$ws.Range("I132").Value = 1243
$ws.Range("J132").Value = 3591.4443
$ws.Range("K132").Value = 3729
$ws.Range("L132").Value = 10774.3329
$ws.Range("M132").Value = -1199
$ws.Range("N132").Value = -15834.3329

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 26749.375
$ws.Range("I61").Value = 30299.428
$ws.Range("K61").Value = 30299.428
$ws.Range("M61").Value = -30097.428

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 26749.375
$ws.Range("I113").Value = 30299.428
$ws.Range("K113").Value = 30299.428
$ws.Range("M113").Value = -28129.428

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3096.1135
$ws.Range("I132").Value = 2441.1353
$ws.Range("K132").Value = 7323.4059
$ws.Range("M132").Value = -4793.4059

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1786.4098
$ws.Range("I136").Value = 1539.5306
$ws.Range("J136").Value = 2794.5
$ws.Range("K136").Value = 4618.5918
$ws.Range("L136").Value = 8383.5
$ws.Range("M136").Value = -2068.5918
$ws.Range("N136").Value = -13483.5

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1602.6428
$ws.Range("I136").Value = 1418.3077
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 4254.9231
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -1704.9231
$ws.Range("N136").Value = -17097

Write-Output "All updates applied"
